$wb = $excel.ActiveWorkbook

# --- Sheet1 (Typography): clear the Widget Wildcard Characters for the "Large" row (H5) ---
$wsTypography = $wb.Worksheets.Item(1)
$wsTypography.Cells.Item(5, 8).ClearContents()

# --- Sheet2 (Translation): update/add rows for the new number keyboard with live view ---
$wsTranslation = $wb.Worksheets.Item(2)

$wsTranslation.Cells.Item(5, 2).Value = "SingleUseId7"
$wsTranslation.Cells.Item(5, 3).Value = "Default"
$wsTranslation.Cells.Item(5, 4).Value = "Center"
$wsTranslation.Cells.Item(5, 5).Value = "LTR"
$wsTranslation.Cells.Item(5, 6).Value = "Witaj w programie inteligenta butelka!`nKlikni START aby rozpocząć <value>"

$wsTranslation.Cells.Item(6, 2).Value = "SingleUseId8"
$wsTranslation.Cells.Item(6, 3).Value = "Default"
$wsTranslation.Cells.Item(6, 4).Value = "Left"
$wsTranslation.Cells.Item(6, 5).Value = "LTR"
$wsTranslation.Cells.Item(6, 6).Value = "Wprowadź dane:"

$wsTranslation.Cells.Item(7, 2).Value = "SingleUseId9"
$wsTranslation.Cells.Item(7, 3).Value = "Default"
$wsTranslation.Cells.Item(7, 4).Value = "Center"
$wsTranslation.Cells.Item(7, 5).Value = "LTR"
$wsTranslation.Cells.Item(7, 6).Value = "Ok"

$wsTranslation.Cells.Item(8, 2).Value = "SingleUseId10"
$wsTranslation.Cells.Item(8, 3).Value = "Default"
$wsTranslation.Cells.Item(8, 4).Value = "Center"
$wsTranslation.Cells.Item(8, 5).Value = "LTR"
$wsTranslation.Cells.Item(8, 6).Value = "1"

$wsTranslation.Cells.Item(9, 2).Value = "SingleUseId11"
$wsTranslation.Cells.Item(9, 3).Value = "Default"
$wsTranslation.Cells.Item(9, 4).Value = "Center"
$wsTranslation.Cells.Item(9, 5).Value = "LTR"
$wsTranslation.Cells.Item(9, 6).Value = "2"

$wsTranslation.Cells.Item(10, 2).Value = "SingleUseId12"
$wsTranslation.Cells.Item(10, 3).Value = "Default"
$wsTranslation.Cells.Item(10, 4).Value = "Center"
$wsTranslation.Cells.Item(10, 5).Value = "LTR"
$wsTranslation.Cells.Item(10, 6).Value = "3"

$wsTranslation.Cells.Item(11, 2).Value = "SingleUseId13"
$wsTranslation.Cells.Item(11, 3).Value = "Default"
$wsTranslation.Cells.Item(11, 4).Value = "Center"
$wsTranslation.Cells.Item(11, 5).Value = "LTR"
$wsTranslation.Cells.Item(11, 6).Value = "4"

$wsTranslation.Cells.Item(12, 2).Value = "SingleUseId14"
$wsTranslation.Cells.Item(12, 3).Value = "Default"
$wsTranslation.Cells.Item(12, 4).Value = "Center"
$wsTranslation.Cells.Item(12, 5).Value = "LTR"
$wsTranslation.Cells.Item(12, 6).Value = "5"

$wsTranslation.Cells.Item(13, 2).Value = "SingleUseId15"
$wsTranslation.Cells.Item(13, 3).Value = "Default"
$wsTranslation.Cells.Item(13, 4).Value = "Center"
$wsTranslation.Cells.Item(13, 5).Value = "LTR"
$wsTranslation.Cells.Item(13, 6).Value = "6"

$wsTranslation.Cells.Item(14, 2).Value = "SingleUseId16"
$wsTranslation.Cells.Item(14, 3).Value = "Default"
$wsTranslation.Cells.Item(14, 4).Value = "Center"
$wsTranslation.Cells.Item(14, 5).Value = "LTR"
$wsTranslation.Cells.Item(14, 6).Value = "7"

$wsTranslation.Cells.Item(15, 2).Value = "SingleUseId17"
$wsTranslation.Cells.Item(15, 3).Value = "Default"
$wsTranslation.Cells.Item(15, 4).Value = "Center"
$wsTranslation.Cells.Item(15, 5).Value = "LTR"
$wsTranslation.Cells.Item(15, 6).Value = "8"

$wsTranslation.Cells.Item(16, 2).Value = "SingleUseId18"
$wsTranslation.Cells.Item(16, 3).Value = "Default"
$wsTranslation.Cells.Item(16, 4).Value = "Center"
$wsTranslation.Cells.Item(16, 5).Value = "LTR"
$wsTranslation.Cells.Item(16, 6).Value = "9"

$wsTranslation.Cells.Item(17, 2).Value = "SingleUseId19"
$wsTranslation.Cells.Item(17, 3).Value = "Default"
$wsTranslation.Cells.Item(17, 4).Value = "Center"
$wsTranslation.Cells.Item(17, 5).Value = "LTR"
$wsTranslation.Cells.Item(17, 6).Value = "0"

$wsTranslation.Cells.Item(18, 2).Value = "SingleUseId22"
$wsTranslation.Cells.Item(18, 3).Value = "Default"
$wsTranslation.Cells.Item(18, 4).Value = "Center"
$wsTranslation.Cells.Item(18, 5).Value = "LTR"
$wsTranslation.Cells.Item(18, 6).Value = "Reset"

$wsTranslation.Cells.Item(19, 2).Value = "SingleUseId23"
$wsTranslation.Cells.Item(19, 3).Value = "Default"
$wsTranslation.Cells.Item(19, 4).Value = "Center"
$wsTranslation.Cells.Item(19, 5).Value = "LTR"
$wsTranslation.Cells.Item(19, 6).Value = "Yes"

$wsTranslation.Cells.Item(20, 2).Value = "SingleUseId24"
$wsTranslation.Cells.Item(20, 3).Value = "Default"
$wsTranslation.Cells.Item(20, 4).Value = "Left"
$wsTranslation.Cells.Item(20, 5).Value = "LTR"
$wsTranslation.Cells.Item(20, 6).Value = "Are you sure to reset?"

$wsTranslation.Cells.Item(21, 2).Value = "SingleUseId25"
$wsTranslation.Cells.Item(21, 3).Value = "Default"
$wsTranslation.Cells.Item(21, 4).Value = "Center"
$wsTranslation.Cells.Item(21, 5).Value = "LTR"
$wsTranslation.Cells.Item(21, 6).Value = "No"

$wsTranslation.Cells.Item(22, 2).Value = "SingleUseId27"
$wsTranslation.Cells.Item(22, 3).Value = "Default"
$wsTranslation.Cells.Item(22, 4).Value = "Left"
$wsTranslation.Cells.Item(22, 5).Value = "LTR"
$wsTranslation.Cells.Item(22, 6).Value = "<value>"

$wsTranslation.Cells.Item(23, 2).Value = "SingleUseId28"
$wsTranslation.Cells.Item(23, 3).Value = "Default"
$wsTranslation.Cells.Item(23, 4).Value = "Left"
$wsTranslation.Cells.Item(23, 5).Value = "LTR"
$wsTranslation.Cells.Item(23, 6).Value = "0"

$wsTranslation.Cells.Item(24, 2).Value = "SingleUseId30"
$wsTranslation.Cells.Item(24, 3).Value = "Default"
$wsTranslation.Cells.Item(24, 4).Value = "Left"
$wsTranslation.Cells.Item(24, 5).Value = "LTR"
$wsTranslation.Cells.Item(24, 6).Value = "<value> cm"

$wsTranslation.Cells.Item(25, 2).Value = "SingleUseId31"
$wsTranslation.Cells.Item(25, 3).Value = "Default"
$wsTranslation.Cells.Item(25, 4).Value = "Left"
$wsTranslation.Cells.Item(25, 5).Value = "LTR"
$wsTranslation.Cells.Item(25, 6).Value = "0"

$wsTranslation.Cells.Item(26, 2).Value = "SingleUseId33"
$wsTranslation.Cells.Item(26, 3).Value = "Default"
$wsTranslation.Cells.Item(26, 4).Value = "Left"
$wsTranslation.Cells.Item(26, 5).Value = "LTR"
$wsTranslation.Cells.Item(26, 6).Value = "<value> l"

$wsTranslation.Cells.Item(27, 2).Value = "SingleUseId34"
$wsTranslation.Cells.Item(27, 3).Value = "Default"
$wsTranslation.Cells.Item(27, 4).Value = "Left"
$wsTranslation.Cells.Item(27, 5).Value = "LTR"
$wsTranslation.Cells.Item(27, 6).Value = "0"

$wsTranslation.Cells.Item(28, 2).Value = "SingleUseId35"
$wsTranslation.Cells.Item(28, 3).Value = "Default"
$wsTranslation.Cells.Item(28, 4).Value = "Center"
$wsTranslation.Cells.Item(28, 5).Value = "LTR"
$wsTranslation.Cells.Item(28, 6).Value = "Waga"

$wsTranslation.Cells.Item(29, 2).Value = "SingleUseId36"
$wsTranslation.Cells.Item(29, 3).Value = "Default"
$wsTranslation.Cells.Item(29, 4).Value = "Center"
$wsTranslation.Cells.Item(29, 5).Value = "LTR"
$wsTranslation.Cells.Item(29, 6).Value = "Wzrost"

$wsTranslation.Cells.Item(30, 2).Value = "SingleUseId37"
$wsTranslation.Cells.Item(30, 3).Value = "Default"
$wsTranslation.Cells.Item(30, 4).Value = "Center"
$wsTranslation.Cells.Item(30, 5).Value = "LTR"
$wsTranslation.Cells.Item(30, 6).Value = "Pojemnosc`n butelki"

$wsTranslation.Cells.Item(31, 2).Value = "SingleUseId38"
$wsTranslation.Cells.Item(31, 3).Value = "Default"
$wsTranslation.Cells.Item(31, 4).Value = "Left"
$wsTranslation.Cells.Item(31, 5).Value = "LTR"
$wsTranslation.Cells.Item(31, 6).Value = "kg"
